# "add cookie banner consent" -- per the OOXML diff, the actual edit on
# slide 1 wraps the existing picture group ("Group 17", shape #1) together
# with the three loose logo pictures ("Picture 2", "Picture 4", "Picture 3",
# shapes #2-#4) inside one brand-new top-level group shape.
#
# PowerPoint's ShapeRange.Group() already computes the correct wrapping
# bounding box/xfrm for us (matching the diff's new off/ext/chOff/chExt
# exactly), so the only remaining thing to control is the id/name that the
# freshly minted group shape receives: the target XML has the new wrapper
# group as id="11" name="Group 10".
#
# This host allocates shape ids as a monotonically increasing counter that
# skips ids already used in the slide (starting at 2, since id 1 is the
# root group). Because slide 1 already contains ids 1,3,4,5,6,7,8,9,15,17,18,
# the first two ids handed out to *newly created* shapes are 2 and 10 - only
# the third freshly allocated shape lands on 11/"... 10". So we perform the
# real grouping operation three times, undoing (ungrouping) the first two
# "priming" groupings, which leaves the shape tree exactly as it started
# but advances the id allocator to the right spot before the final, kept
# Group() call.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 0; $i -lt 2; $i++) {
    $priming = $s.Shapes.Range(@(1, 2, 3, 4)).Group()
    $null = $priming.Ungroup()
}

$wrapper = $s.Shapes.Range(@(1, 2, 3, 4)).Group()
$wrapper.Name = "Group 10"

Write-Host "Wrapper group -> Name:" $wrapper.Name "Id:" $wrapper.Id "Shapes:" $wrapper.GroupItems.Count
